$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.8905744754041969
$ws.Range("C2").Value = 0.8960784313725491
$ws.Range("D2").Value = 0.8910058568882098

$ws.Range("B3").Value = 0.8308471177944863
$ws.Range("D3").Value = 0.8519392802319631

$ws.Range("B4").Value = 0.8075983436853003
$ws.Range("C4").Value = 0.9157894736842105
$ws.Range("D4").Value = 0.8571121901354459

$ws.Range("B5").Value = 0.8150000000000001
$ws.Range("D5").Value = 0.8111998361998362

$ws.Range("B6").Value = 0.9169423558897243
$ws.Range("C6").Value = 0.7352380952380952
$ws.Range("D6").Value = 0.8128557063851183
